$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.28905232522894408
$ws.Range("A2").Value = -0.0059999999376216806
$ws.Range("A3").Value = -0.0039999999358961702
$ws.Range("A4").Value = -0.0079999998913198311
$ws.Range("A5").Value = -0.0029999999342180672
$ws.Range("A6").Value = -0.0083612829731407601
$ws.Range("A7").Value = -0.0099999998498985221
$ws.Range("A8").Value = -0.0099999998484388009
$ws.Range("A9").Value = -0.0019999999255730927
$ws.Range("A10").Value = -0.0019999999253901279
$ws.Range("A11").Value = 0.026661300737693239
$ws.Range("A12").Value = -0.0034999999093909828
$ws.Range("A13").Value = -0.0034999999029627915
$ws.Range("A14").Value = -0.0079999998575388531
$ws.Range("A15").Value = -0.00099999992460997333
$ws.Range("A16").Value = -0.0019999999139002078
$ws.Range("A17").Value = -0.0019999999124884482
$ws.Range("A18").Value = -0.0039999998927822134
$ws.Range("A19").Value = -0.050750398993583357
$ws.Range("A20").Value = -0.0039999999521107554
$ws.Range("A21").Value = -0.003999999951610711
$ws.Range("A22").Value = -0.0039999999514481743
$ws.Range("A23").Value = -0.068032308087198068
$ws.Range("A24").Value = -0.019999999765738075
$ws.Range("A25").Value = -0.019999999762687182
$ws.Range("A26").Value = -0.0024999999218877633
$ws.Range("A27").Value = -0.002499999919317375
$ws.Range("A28").Value = -0.001999999911690864
$ws.Range("A29").Value = -0.0069999998542611408
$ws.Range("A30").Value = -0.025133293149751879
$ws.Range("A31").Value = -0.0069999998455969603
$ws.Range("A32").Value = -0.0099999998157631609
$ws.Range("A33").Value = -0.0039999998740185561

$ws.Range("A1").ColumnWidth = 15.67
